# Update input data, rename date column
$wb = $excel.ActiveWorkbook

# "ColumnHeaders" is the first sheet (sheetId=1 / rId1); rename the
# "datetime" attribute row to "date_time_utc".
$ws = $wb.Worksheets.Item("ColumnHeaders")
$ws.Range("A6").Value = "date_time_utc"

# Move the active selection as left by the edit (was A5, now B27).
$ws.Range("B27").Select()
